$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename the IMF / OECD(20%) columns -----------------------
# F/G used to be "IMF - Sales" / "IMF - Sales + Emp" -> become "IMF (20%) - ..."
# H/I used to be "OECD (20%) - Sales" / "OECD (20%) - Sales + Emp" -> become "IMF - ..."
# (J/K "OECD - Sales" / "OECD - Sales + Emp" stay untouched)
$ws.Range("H1").Value = "IMF - Sales"
$ws.Range("I1").Value = "IMF - Sales + Emp"
$ws.Range("F1").Value = "IMF (20%) - Sales"
$ws.Range("G1").Value = "IMF (20%) - Sales + Emp"

# --- Data rows 2-11 ---------------------------------------------------------
# The former F/G values (old "IMF - Sales(+Emp)" data) shift right into H/I.
# New F/G values are freshly computed "IMF (20%)" figures.
# The former H/I values (old "OECD (20%)" data) are discarded entirely.

$oldF = @{
    2  = 1.830532502618961
    3  = 2.749696733824161
    4  = 0.505747365869219
    5  = 0.8574946411960901
    6  = 4.250911976949065
    7  = 1.329521722860183
    8  = 2.061815722023631
    9  = 3.186771358099338
    10 = 1.465451043466264
    11 = 7.431748429398017
}

$oldG = @{
    2  = 1.497187058396667
    3  = 3.147379767852597
    4  = 2.405701843541175
    5  = 3.543764499696386
    6  = 3.228815446321175
    7  = 6.700190657546217
    8  = 2.268948420676758
    9  = 2.647682377484735
    10 = 1.846481599431316
    11 = 4.623537134690935
}

$newF = @{
    2  = 0.366106500523793
    3  = 0.5499393467648322
    4  = 0.1011494731738439
    5  = 0.1714989282392181
    6  = 0.8501823953898133
    7  = 0.2659043445720369
    8  = 0.4123631444047267
    9  = 0.6373542716198679
    10 = 0.2930902086932531
    11 = 1.486349685879603
}

$newG = @{
    2  = 0.2994374116793341
    3  = 0.629475953570518
    4  = 0.4811403687082352
    5  = 0.7087528999392775
    6  = 0.6457630892642331
    7  = 1.340038131509243
    8  = 0.4537896841353516
    9  = 0.5295364754969457
    10 = 0.3692963198862632
    11 = 0.9247074269381814
}

for ($r = 2; $r -le 11; $r++) {
    $ws.Range("H$r").Value = $oldF[$r]
    $ws.Range("I$r").Value = $oldG[$r]
    $ws.Range("F$r").Value = $newF[$r]
    $ws.Range("G$r").Value = $newG[$r]
}
